$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 10:31"

# --- Daily refresh of case numbers (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---------------------

# India (row 6)
$ws.Range("B6").Value = 2091416
$ws.Range("C6").Value = 4552
$ws.Range("E6").Value = 619699
$ws.Range("G6").Value = 39
$ws.Range("H6").Value = 42617

# Rusia (row 7)
$ws.Range("B7").Value = 882347
$ws.Range("C7").Value = 5212
$ws.Range("D7").Value = 690207
$ws.Range("E7").Value = 177286
$ws.Range("G7").Value = 129
$ws.Range("H7").Value = 14854

# Singapur (row 47)
$ws.Range("B47").Value = 54929
$ws.Range("C47").Value = 132
$ws.Range("E47").Value = 6590

# Armenia overtakes Ghana -> row 54 becomes Armenia (updated numbers),
# row 55 becomes Ghana (the numbers Armenia's old row held before the swap
# settle to what used to be row 54's figures).
$ws.Range("A54").Value = "Armenia"
$ws.Range("B54").Value = 40185
$ws.Range("C54").Value = 200
$ws.Range("D54").Value = 32395
$ws.Range("E54").Value = 7005
$ws.Range("G54").Value = 8
$ws.Range("H54").Value = 785

$ws.Range("A55").Value = "Ghana"
$ws.Range("B55").Value = 40097
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 36638
$ws.Range("E55").Value = 3253
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 206

# Austria (row 70)
$ws.Range("B70").Value = 21919
$ws.Range("C70").Value = 82
$ws.Range("D70").Value = 19812
$ws.Range("E70").Value = 1386
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 721

# Hungria overtakes Republica de Africa Central -> row 107 becomes Hungria
# (updated numbers), row 108 becomes Republica de Africa Central (the
# figures that used to sit in row 107).
$ws.Range("A107").Value = "Hungria"
$ws.Range("B107").Value = 4653
$ws.Range("C107").Value = 32
$ws.Range("D107").Value = 3491
$ws.Range("E107").Value = 560
$ws.Range("H107").Value = 602

$ws.Range("A108").Value = "Republica de Africa Central"
$ws.Range("B108").Value = 4641
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 1716
$ws.Range("E108").Value = 2866
$ws.Range("H108").Value = 59

# Hong Kong (row 111)
$ws.Range("E111").Value = 1272
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 47

# Santa Lucia overtakes Timor Oriental (identical figures, just the name
# order swaps: row 202 becomes Santa Lucia, row 203 becomes Timor Oriental).
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"
